# Daily attendance processing - 2026-01-25 09:34:59
#
# The "Recorded By" column (G) stores a comma-separated list such as
# "System, dnasr281@gmail.com". For every row where that list is exactly
# "System, dnasr281@gmail.com", flip it to "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
